$wb = $excel.ActiveWorkbook

# Update the "Ready for handoff" status text to "In Translation" everywhere
# it appears: the Overview sheet's per-language status columns, and the
# "Status" column on each per-language detail sheet.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# The shortened text means the status columns no longer need to be as wide;
# re-fit them to the new content.
$wsOverview.Range("E:F").Columns.AutoFit()
$wsZhCn.Range("C:C").Columns.AutoFit()
$wsDeDe.Range("C:C").Columns.AutoFit()
